$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.911.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.792.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5125"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.02%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.94"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.243"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9998"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.239"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.779.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001076"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06525"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.932"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.980.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.229"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.04"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.993.73"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.373"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1078"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.041"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.610"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.492"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07088"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02308"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.747"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2130"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.020"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6099"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9992"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.310"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5918"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.700"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.911"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06822"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.41%  "
